$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.891.84'
$ws.Range("E2").Value = '  -0.25%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.633.74'
$ws.Range("E3").Value = '  -0.51%  '

$ws.Range("E4").Value = '  -0.25%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '216.06'
$ws.Range("E5").Value = '  +0.35%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5074'
$ws.Range("E6").Value = '  +0.17%  '

$ws.Range("E7").Value = '  -0.12%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2575'
$ws.Range("E8").Value = '  +0.69%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06360'
$ws.Range("E9").Value = '  -0.21%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.59'
$ws.Range("E10").Value = '  +0.65%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07759'
$ws.Range("E11").Value = '  -0.07%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.253'
$ws.Range("E12").Value = '  -0.60%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.641.21'
$ws.Range("E13").Value = '  -0.76%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5503'
$ws.Range("E14").Value = '  +0.81%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0₅7700'
$ws.Range("E15").Value = '  -1.59%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.83'
$ws.Range("E16").Value = '  -0.70%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '25.916.88'
$ws.Range("E17").Value = '  -0.29%  '

$ws.Range("E18").Value = '  -0.02%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.432'
$ws.Range("E19").Value = '  -0.27%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '194.34'
$ws.Range("E20").Value = '  -1.53%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.908'
$ws.Range("E21").Value = '  -0.34%  '

$ws.Range("E22").Value = '  +0.04%  '

$ws.Range("E23").Value = '  -0.24%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.904'
$ws.Range("E24").Value = '  +0.57%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '142.58'
$ws.Range("E25").Value = '  +1.05%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1239'
$ws.Range("E26").Value = '  +6.06%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.807'
$ws.Range("E27").Value = '  -1.16%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.59'
$ws.Range("E28").Value = '  -0.60%  '

$ws.Range("E29").Value = '  +0.27%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.04892'
$ws.Range("E30").Value = '  -1.28%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.254'
$ws.Range("E31").Value = '  -0.16%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.188'
$ws.Range("E32").Value = '  +0.14%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.542'
$ws.Range("E33").Value = '  +0.17%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.381'
$ws.Range("E34").Value = '  +0.71%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9026'
$ws.Range("E35").Value = '  +0.83%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.559'
$ws.Range("E36").Value = '  -1.18%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("B37").Value = 'Maker'
$ws.Range("C37").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D37").Value = '1.125.44'
$ws.Range("E37").Value = '  -0.75%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("B38").Value = 'ImmutableX'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D38").Value = '0.5502'
$ws.Range("E38").Value = '  +1.20%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01555'
$ws.Range("E39").Value = '  -0.30%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.002'
$ws.Range("E40").Value = '  -0.15%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.586'
$ws.Range("E41").Value = '  +0.01%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8032'
$ws.Range("E42").Value = '  -2.18%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '97.56'
$ws.Range("E43").Value = '  -1.98%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0₈120'
$ws.Range("E44").Value = '  -5.30%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.772.68'
$ws.Range("E45").Value = '  -0.28%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4452'
$ws.Range("E46").Value = '  -1.98%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '54.77'
$ws.Range("E47").Value = '  -0.01%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.9950'
$ws.Range("E48").Value = '  -0.87%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05146'
$ws.Range("E49").Value = '  +1.43%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.530'
$ws.Range("E50").Value = '  +1.71%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.002'
$ws.Range("E51").Value = '  -0.35%  '
